$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '22.027.96'
$ws.Range("E2").Value = '  -0.37%  '
$ws.Range("D3").Value = '1.553.89'
$ws.Range("E3").Value = '  +0.19%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.28%  '
$ws.Range("E5").Value = '  -0.11%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '289.99'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.79%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3967'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +4.11%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3213'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.88%  '
$ws.Range("E9").Value = '  -0.07%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07220'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.81%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.075'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -5.00%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.000'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.30%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.680'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.28%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '18.65'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -7.02%  '
$ws.Range("E15").Value = '  +5.35%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.614'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.56%  '
$ws.Range("D17").Value = '1.552.09'
$ws.Range("E17").Value = '  -2.10%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06585'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.95%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '83.42'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.70%  '
$ws.Range("E20").Value = '  -0.25%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.245'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.80%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '15.48'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.48%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.29'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.26%  '
$ws.Range("D24").Value = '22.038.64'
$ws.Range("E24").Value = '  -0.29%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.362'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.26%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.410'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.69%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '148.64'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.41%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.52'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.24%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.880'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.63%  '
$ws.Range("D30").Value = '1.726.33'
$ws.Range("E30").Value = '  -1.61%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '118.78'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.77%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9785'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -10.24%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.792'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.10%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08320'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.60%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.181'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.52%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.600'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -15.78%  '
$ws.Range("E37").Value = '  -2.81%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.087'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.81%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05980'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.90%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.212'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.32%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.2032'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.31%  '
$ws.Range("E42").Value = '  -0.17%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '10.70'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.07%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5804'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.91%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.739'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.16%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.93'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -6.45%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5554'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.84%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '117.84'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.06%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.892'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.53%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.133'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.23%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06815'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.94%  '
